$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.486941814422607
$ws.Range("B1").Value = 3.583910465240479
$ws.Range("C1").Value = 2.564816236495972
$ws.Range("D1").Value = 1.322996735572815
$ws.Range("E1").Value = 0.7646856904029846
